# Update "想去人数" (attendance count) figures on the 展览 and 全部类型 sheets
# to match the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 874
$ws1.Range("F3").Value = 1449
$ws1.Range("F4").Value = 1105
$ws1.Range("F5").Value = 523
$ws1.Range("F14").Value = 2286
$ws1.Range("F18").Value = 276
$ws1.Range("F24").Value = 244
$ws1.Range("F27").Value = 1571
$ws1.Range("F28").Value = 307

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 874
$ws4.Range("F4").Value = 1449
$ws4.Range("F5").Value = 1105
$ws4.Range("F8").Value = 523
$ws4.Range("F18").Value = 2286
$ws4.Range("F23").Value = 276
$ws4.Range("F37").Value = 244
$ws4.Range("F40").Value = 1571
$ws4.Range("F41").Value = 307
